$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 4 currently reads "try: n > 0 :" (a leftover coding mistake).
# Keep the "try: " portion and remove the erroneous "n > 0 :" tail so the
# paragraph reads just "try: ".
$para = $tr.Paragraphs(4, 1)
$bad = $para.Characters(6, $para.Length - 5)
$bad.Text = ""
